# Update the monitoring data: re-rank empadronadores by total_registros
# and refresh the counts according to the new data upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ALBIRENA GARCIA ANGEELO ALONSO", 187),
    @("URRIOLA ARISMENDIZ INGRID MARYURI", 162),
    @("MARYURI OJEDA VALLE", 160),
    @("ALAMA NIMA CLARITZA MABEL", 157),
    @("CORDOVA CARMEN ANGIE NATALLY", 156),
    @("ATOCHE PALACIOS LUIS ANGEL", 155),
    @("AGURTO ORDINOLA LISBET JAQUELIN", 149),
    @("MANUEL LEUNARDO PRADO BAILON", 147),
    @("RUIDIAS FRIAS MELISSA VICTORIA", 139),
    @("VEGA ROBLEDO FERNANDO ERNESTO", 136),
    @("JUAREZ CARMEN PIERRE ALEXANDER", 135),
    @("ROMAN GALECIO MARITZA DEL PILAR", 134),
    @("BERNAOLA CARMEN ZUMIKO YASHURY", 119),
    @("CARREÑO PALACIOS KATHERINE DE LOS MILAGROS", 100)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
